# Apply weekly update: rotate the variable fields (D, L, M, N, O, P, Q, S, T)
# across rows 2-9 so that each row now shows the data that previously
# belonged to the row 5 positions later in the 8-row cycle (2..9).
#
# This reproduces the diff where row 2 takes what was row 7's data, row 3
# takes what was row 8's, ..., row 6 takes what was row 3's data, etc.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot of the "before" values for the fields that change, keyed by row.
$data = @{
    2 = @{ D = 44230; L = "Primera";  M = 160; N = 16500; O = 17000; P = 16750; Q = "$/caja 18 kilos granel"; S = 931;  T = 18 }
    3 = @{ D = 44230; L = "Segunda";  M = 160; N = 14500; O = 15000; P = 14750; Q = "$/caja 18 kilos granel"; S = 819;  T = 18 }
    4 = @{ D = 44224; L = "Especial"; M = 100; N = 16500; O = 17000; P = 16750; Q = "$/caja 16 kilos granel"; S = 1047; T = 16 }
    5 = @{ D = 44224; L = "Primera";  M = 200; N = 14500; O = 15000; P = 14750; Q = "$/caja 16 kilos granel"; S = 922;  T = 16 }
    6 = @{ D = 44224; L = "Segunda";  M = 200; N = 12500; O = 13000; P = 12750; Q = "$/caja 16 kilos granel"; S = 797;  T = 16 }
    7 = @{ D = 44209; L = "Primera";  M = 300; N = 15500; O = 16000; P = 15750; Q = "$/caja 16 kilos granel"; S = 984;  T = 16 }
    8 = @{ D = 44210; L = "Primera";  M = 240; N = 15500; O = 16000; P = 15750; Q = "$/caja 16 kilos granel"; S = 984;  T = 16 }
    9 = @{ D = 44210; L = "Segunda";  M = 300; N = 12500; O = 13000; P = 12750; Q = "$/caja 16 kilos granel"; S = 797;  T = 16 }
}

# Mapping from target row to the source row whose data it should now hold.
$mapping = @{ 2 = 7; 3 = 8; 4 = 9; 5 = 2; 6 = 3; 7 = 4; 8 = 5; 9 = 6 }

foreach ($row in 2..9) {
    $src = $data[$mapping[$row]]

    $ws.Cells.Item($row, 4).Value2 = $src.D   # D: Fecha
    $ws.Range("L$row").Value = $src.L        # L: Calidad
    $ws.Range("M$row").Value = $src.M        # M: Volumen
    $ws.Range("N$row").Value = $src.N        # N: Precio minimo
    $ws.Range("O$row").Value = $src.O        # O: Precio maximo
    $ws.Range("P$row").Value = $src.P        # P: Precio promedio ponderado
    $ws.Range("Q$row").Value = $src.Q        # Q: Unidad de comercializacion
    $ws.Range("S$row").Value = $src.S        # S: Precio $/Kg
    $ws.Range("T$row").Value = $src.T        # T: Kg / unidad
}
